$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - "Making postman collection" entry
$ws.Range("D7").Value = "Today I wanted to work on the postman collection because this makes it easier to test. The routes are all explained as I thought needed. Some have more explanation than others but this is due to the difficulty of the route."
$ws.Range("D7").WrapText = $true
$ws.Range("A7").Value = "Making postman collection"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 44739
$ws.Range("E7").Value = "x"

# Row 8 - "Implementing cursor pagination" entry
$ws.Range("A8").Value = "Implementing cursor pagination"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 44740
$ws.Range("D8").Value = "I am not entirely sure about this but I think that I have implemented the cursor pagination. I have never implemented cursor pagination before so this was very new to me. First I started by researching what cursor pagination is. I followed many different courses, videos and documentations. After a few minutes I got the feeling that the concept was pretty clear to me. At that moment I started programming and searching online for different kinds of examples. "
$ws.Range("D8").WrapText = $true
$ws.Range("E8").Value = "x"

# Update selection / view state (scroll so row 2 is the top visible row, like the source workbook)
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D13").Select()
